$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm2"
$ws.Range("C2").Value = "Ramp3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.072366
$ws.Range("H2").Value = 0.217098
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2287846666666667
$ws.Range("N2").Value = 0.686354
$ws.Range("O2").Value = 0.2689534993532379
$ws.Range("P2").Value = 0.268953499353238
$ws.Range("Q2").Value = 0.016556231188
$ws.Range("R2").Value = 0.149006080692
$ws.Range("S2").Value = 0.2689534993532379
$ws.Range("T2").Value = 0.268953499353238

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adm2"
$ws.Range("C3").Value = "Ramp3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.072366
$ws.Range("H3").Value = 0.217098
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.455408
$ws.Range("N3").Value = 1.366224
$ws.Range("O3").Value = 0.5353661896053321
$ws.Range("P3").Value = 0.5353661896053322
$ws.Range("Q3").Value = 0.032956055328
$ws.Range("R3").Value = 0.2966044979520001
$ws.Range("S3").Value = 0.5353661896053321
$ws.Range("T3").Value = 0.5353661896053322

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adm2"
$ws.Range("C4").Value = "Ramp3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.072366
$ws.Range("H4").Value = 0.217098
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04315
$ws.Range("N4").Value = 0.12945
$ws.Range("O4").Value = 0.05072605461799107
$ws.Range("P4").Value = 0.05072605461799109
$ws.Range("Q4").Value = 0.0031225929
$ws.Range("R4").Value = 0.0281033361
$ws.Range("S4").Value = 0.05072605461799107
$ws.Range("T4").Value = 0.05072605461799109

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Adm2"
$ws.Range("C5").Value = "Ramp3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.072366
$ws.Range("H5").Value = 0.217098
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.123305
$ws.Range("N5").Value = 0.369915
$ws.Range("O5").Value = 0.1449542564234389
$ws.Range("P5").Value = 0.144954256423439
$ws.Range("Q5").Value = 0.00892308963
$ws.Range("R5").Value = 0.08030780667000001
$ws.Range("S5").Value = 0.1449542564234389
$ws.Range("T5").Value = 0.144954256423439
